# Add a new "Эксплуатация объекта" measure row (id 600) as row 7 on the
# first worksheet, mirroring the existing id/name rows above it, then
# update the page setup to match an A4 portrait layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row: A7 = 600, B7 = "Эксплуатация объекта" (appended as a new
# shared string, same pattern as the existing rows 2-6).
$ws.Range("A7").Value = 600
$ws.Range("B7").Value = "Эксплуатация объекта"

# Move/collapse the selection onto the newly added cell, as Excel does
# after editing the last row of a table.
$ws.Range("B7").Select()

# Set printer/page setup (A4, portrait) for the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
